# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1) Metadata sheet: bump the "Date" property value.
# 2) Elements sheet: append a new mapping column
#    "Mapping: Spécification métier vers l'extension AsLieuDit"
#    with values for the 5 data rows (only the last row - the
#    Extension.value[x] element - gets a mapping: "lieuDit").

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date property -----------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: add the new mapping column ----------------------------
$ws = $wb.Worksheets.Item("Elements")

# Clone column AK (the last existing "Mapping: ..." column) into the new
# column AL so the new column inherits the correct header/body cell
# formatting (style "1" for the header row, style "2" for data rows).
$ws.Range("AK1:AK6").Copy($ws.Range("AL1:AL6"))

# Now overwrite the copied values with the real content for the new column.
$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension AsLieuDit"
$ws.Range("AL2").Value = ""
$ws.Range("AL3").Value = ""
$ws.Range("AL4").Value = ""
$ws.Range("AL5").Value = ""
$ws.Range("AL6").Value = "lieuDit"

# Match the column's width from the source workbook (as closely as the
# engine's character-width rounding allows).
$ws.Columns.Item(38).ColumnWidth = 59
